# Update the "Förändrad" (Changed) date column (C) for every data row.
# The workbook records when each row was last refreshed; this run refreshed
# the data on a later date, so column C moves from 2023-09-03 (45172) to
# 2023-09-06 (45175) for all data rows (2 through 490).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 490 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45172) {
        $cell.Value2 = 45175
    }
}
